$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top 10 Cities Least Debt")

# Insert a new row at row 3, shifting rows 3-11 down to 4-12
$ws.Rows.Item(3).Insert()

# Set the new row 3 values (spokane, Washington)
$ws.Cells.Item(3, 1).Value = "spokane"
$ws.Cells.Item(3, 2).Value = "Washington"
$ws.Cells.Item(3, 3).Value = -196311784
$ws.Cells.Item(3, 4).Value = -857.3390631414371

# Remove the last row (previously row 11 "aurora", now shifted to row 12) so only 10 data rows remain
$ws.Rows.Item(12).Delete()
